$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 keeps referencing the pre-existing "Title was choosen" string but is re-set here too.
# Order below matches the order new unique strings were first introduced by the author
# so the shared-strings table comes out in the same sequence.

# Row 3: B3 -> "Concepts of login page.", C3 -> "Feasibility study on this project"
$ws.Range("B3").Value = "Concepts of login page."
$ws.Range("C3").Value = "Feasibility study on this project"

# Row 2: B2 -> "Read concepts..." (new wording), C2 -> "Title was choosen"
$ws.Range("B2").Value = "Read concepts of how to create a login page, went through the abstract of previous projects."
$ws.Range("C2").Value = "Title was choosen"

# Row 3 (cont.): D3 -> new text
$ws.Range("D3").Value = "If password is forget during login, then OTP is sent to the respective mail id"

# Row 4: new row
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 43808
$ws.Range("B4").Value = "Feasibity study done. Prototype for the project was developed"
$ws.Range("C4").Value = "Modules, prototype was done"

# Column width adjustments
$ws.Range("A1").EntireColumn.ColumnWidth = 11.666666666666666
$ws.Range("D1").EntireColumn.ColumnWidth = 67.16666666666667

# View adjustments
$excel.ActiveWindow.Zoom = 113
[void]$ws.Range("C4").Select()
